$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data row 42 (2025, Q1) for WZ08-2221(2/3) + WZ08-2222 / IK Kunststoffverpackungen
$ws.Range("A42").Value = 2025
$ws.Range("B42").Value = "Q1"
$ws.Range("G42").Value = "WZ08-2221(2/3) + WZ08-2222"
$ws.Range("H42").Value = "IK Kunststoffverpackungen"
$ws.Range("I42").Value = -17.39130434782609
$ws.Range("J42").Value = -43.478260869565219
$ws.Range("K42").Value = -9.8901098901098905
$ws.Range("L42").Value = 5.4347826086956523
$ws.Range("M42").Value = -7.608695652173914
$ws.Range("N42").Value = 1.086956521739129
$ws.Range("O42").Value = -66.304347826086953

# Update view state: selection/top-left cell to match target workbook
$ws.Range("X45").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 4
